$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'27.286.13"
$ws.Range('E2').Value = "'  -2.40%  "

# Row 3
$ws.Range('D3').Value = "'1.706.62"
$ws.Range('E3').Value = "'  -1.82%  "

# Row 4
$ws.Range('E4').Value = "'  +0.03%  "

# Row 5
$ws.Range('D5').Value = "'223.35"
$ws.Range('E5').Value = "'  -2.70%  "

# Row 6
$ws.Range('D6').Value = "'0.5305"
$ws.Range('E6').Value = "'  -2.58%  "

# Row 7
$ws.Range('D7').Value = "'1.004"
$ws.Range('E7').Value = "'  +0.14%  "

# Row 8
$ws.Range('D8').Value = "'0.2651"
$ws.Range('E8').Value = "'  -4.66%  "

# Row 9
$ws.Range('D9').Value = "'0.06584"
$ws.Range('E9').Value = "'  -1.92%  "

# Row 10
$ws.Range('D10').Value = "'20.84"
$ws.Range('E10').Value = "'  -4.19%  "

# Row 11
$ws.Range('D11').Value = "'0.07627"
$ws.Range('E11').Value = "'  -2.02%  "

# Row 12
$ws.Range('D12').Value = "'4.576"
$ws.Range('E12').Value = "'  -2.91%  "

# Row 13
$ws.Range('B13').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C13').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D13').Value = "'1.944.49"
$ws.Range('E13').Value = "'  -1.70%  "

# Row 14
$ws.Range('B14').Value = "'WrappedEther"
$ws.Range('C14').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D14').Value = "'1.678.23"
$ws.Range('E14').Value = "'  -3.44%  "

# Row 15
$ws.Range('D15').Value = "'0.5728"
$ws.Range('E15').Value = "'  -4.51%  "

# Row 16
$ws.Range('D16').Value = "'0.0₅8179"
$ws.Range('E16').Value = "'  -2.97%  "

# Row 17
$ws.Range('D17').Value = "'67.53"
$ws.Range('E17').Value = "'  -3.12%  "

# Row 18
$ws.Range('D18').Value = "'27.283.71"
$ws.Range('E18').Value = "'  -2.33%  "

# Row 19
$ws.Range('D19').Value = "'216.06"
$ws.Range('E19').Value = "'  -3.75%  "

# Row 20
$ws.Range('E20').Value = "'  +0.04%  "

# Row 21
$ws.Range('D21').Value = "'4.667"
$ws.Range('E21').Value = "'  -3.45%  "

# Row 22
$ws.Range('E22').Value = "'  -4.96%  "

# Row 23
$ws.Range('D23').Value = "'5.966"
$ws.Range('E23').Value = "'  -4.57%  "

# Row 24
$ws.Range('E24').Value = "'  +0.00%  "

# Row 25
$ws.Range('D25').Value = "'141.84"
$ws.Range('E25').Value = "'  -3.34%  "

# Row 26
$ws.Range('D26').Value = "'1.756"
$ws.Range('E26').Value = "'  +6.57%  "

# Row 27
$ws.Range('D27').Value = "'0.1216"
$ws.Range('E27').Value = "'  -2.69%  "

# Row 28
$ws.Range('D28').Value = "'7.257"
$ws.Range('E28').Value = "'  -2.72%  "

# Row 29
$ws.Range('D29').Value = "'16.30"
$ws.Range('E29').Value = "'  -4.84%  "

# Row 30
$ws.Range('D30').Value = "'0.05386"
$ws.Range('E30').Value = "'  -4.76%  "

# Row 31
$ws.Range('D31').Value = "'1.292"
$ws.Range('E31').Value = "'  -2.10%  "

# Row 32
$ws.Range('D32').Value = "'3.501"
$ws.Range('E32').Value = "'  -5.36%  "

# Row 33
$ws.Range('D33').Value = "'3.419"
$ws.Range('E33').Value = "'  -3.69%  "

# Row 34
$ws.Range('D34').Value = "'1.639"
$ws.Range('E34').Value = "'  -1.48%  "

# Row 35
$ws.Range('D35').Value = "'2.879"
$ws.Range('E35').Value = "'  +0.64%  "

# Row 36
$ws.Range('D36').Value = "'2.424"
$ws.Range('E36').Value = "'  -1.14%  "

# Row 37
$ws.Range('D37').Value = "'0.9463"
$ws.Range('E37').Value = "'  -3.99%  "

# Row 38
$ws.Range('D38').Value = "'0.5843"
$ws.Range('E38').Value = "'  -1.85%  "

# Row 39
$ws.Range('D39').Value = "'0.01633"
$ws.Range('E39').Value = "'  -2.46%  "

# Row 40
$ws.Range('E40').Value = "'  -2.75%  "

# Row 41
$ws.Range('B41').Value = "'PaxDollar"
$ws.Range('C41').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('D41').Value = "'1.004"
$ws.Range('E41').Value = "'  +0.02%  "

# Row 42
$ws.Range('B42').Value = "'Maker"
$ws.Range('C42').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D42').Value = "'1.042.75"
$ws.Range('E42').Value = "'  -0.49%  "

# Row 43
$ws.Range('D43').Value = "'0.8388"
$ws.Range('E43').Value = "'  -1.08%  "

# Row 44
$ws.Range('D44').Value = "'100.97"
$ws.Range('E44').Value = "'  -1.34%  "

# Row 45
$ws.Range('D45').Value = "'1.849.98"
$ws.Range('E45').Value = "'  -1.77%  "

# Row 46
$ws.Range('E46').Value = "'  +4.41%  "

# Row 47
$ws.Range('D47').Value = "'58.02"
$ws.Range('E47').Value = "'  -3.61%  "

# Row 48
$ws.Range('D48').Value = "'0.4502"
$ws.Range('E48').Value = "'  +1.68%  "

# Row 49
$ws.Range('E49').Value = "'  -0.80%  "

# Row 50
$ws.Range('D50').Value = "'8.099"
$ws.Range('E50').Value = "'  -2.67%  "

# Row 51
$ws.Range('B51').Value = "'Cronos"
$ws.Range('C51').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D51').Value = "'0.05241"
$ws.Range('E51').Value = "'  -1.45%  "
